$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 278, shifting existing rows 278-399 down to 279-400.
$ws.Rows.Item(278).Insert()

# Populate the newly inserted row 278 with the new data record.
$ws.Cells.Item(278, 1).Value = 10
$ws.Cells.Item(278, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(278, 3).Value = "La Araucanía"
$ws.Cells.Item(278, 4).Value = 44784
$ws.Cells.Item(278, 5).Value = 9
$ws.Cells.Item(278, 6).Value = "Fruta"
$ws.Cells.Item(278, 7).Value = 100108
$ws.Cells.Item(278, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(278, 9).Value = 100108002
$ws.Cells.Item(278, 10).Value = "Mango"
$ws.Cells.Item(278, 11).Value = "Sin especificar"
$ws.Cells.Item(278, 12).Value = "Primera"
$ws.Cells.Item(278, 13).Value = 145
$ws.Cells.Item(278, 14).Value = 10000
$ws.Cells.Item(278, 15).Value = 11000
$ws.Cells.Item(278, 16).Value = 10414
$ws.Cells.Item(278, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(278, 18).Value = "Brasil"
$ws.Cells.Item(278, 19).Value = 2604
$ws.Cells.Item(278, 20).Value = 4
